# Adds new English/Russian word pairs to the vocabulary sheet.
# 1) Remove the extra duplicate "Word"/"Слово" rows left over from a
#    data-entry glitch (rows 13, 14 and 15 were all "Word", right before
#    "Until the end of time" / "Exception" / another "Word" row), leaving
#    a single "Word" row, matching the corrected source data.
# 2) Append the new word pairs at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove two of the duplicated "Word" rows (rows 13 and 14), leaving a
# single "Word" row (originally row 15) just before "Restore".
$ws.Rows("13:15").Delete()

$newWords = @(
    @("Term", "Условие"),
    @("Specification", "Особенность"),
    @("Responsibilities", "Обязанности"),
    @("Requirements", "Требования"),
    @("Reliable", "Надёжный"),
    @("Purpose", "Цель"),
    @("Possible", "Возможный"),
    @("Maintenance", "Техобслуживание"),
    @("Launch", "Запуск"),
    @("Itemize", "Составить перечень"),
    @("Indicate", "Указать"),
    @("Establish", "Устанавливать"),
    @("Edit", "Редактировать"),
    @("Fix", "Исправлять"),
    @("Domain", "Домен"),
    @("Displayed", "Отображаемый"),
    @("Discretion", "Усмотрение"),
    @("Database", "База данных"),
    @("Cover", "Описывать"),
    @("Content", "Контент"),
    @("Compatibility", "Совместимость"),
    @("Brief", "Краткий"),
    @("Carefully", "Тщательно"),
    @("Beta test", "Опытная эксплуатация"),
    @("Backup", "Запас"),
    @("Adjusted", "Настроенный"),
    @("Visual aids", "Наглядные пособия"),
    @("Piace", "Мир"),
    @("Extra randomness", "Лишняя случайность"),
    @("Check out the documentation", "Ознакомиться с документацией"),
    @("Multithreading", "Многопоточность")
)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$row = $lastRow + 1
foreach ($pair in $newWords) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}

